# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets
# to match the newly scraped numbers (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value  = 6645
$wsExhibit.Range("F4").Value  = 418
$wsExhibit.Range("F7").Value  = 542
$wsExhibit.Range("F8").Value  = 98
$wsExhibit.Range("F12").Value = 168
$wsExhibit.Range("F13").Value = 396
$wsExhibit.Range("F14").Value = 1286
$wsExhibit.Range("F16").Value = 3314
$wsExhibit.Range("F19").Value = 1952
$wsExhibit.Range("F20").Value = 69
$wsExhibit.Range("F21").Value = 27
$wsExhibit.Range("F22").Value = 128

# --- Sheet "全部类型" (all types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 6645
$wsAll.Range("F4").Value  = 418
$wsAll.Range("F8").Value  = 542
$wsAll.Range("F9").Value  = 98
$wsAll.Range("F13").Value = 168
$wsAll.Range("F14").Value = 396
$wsAll.Range("F15").Value = 1286
$wsAll.Range("F17").Value = 3314
$wsAll.Range("F20").Value = 1952
$wsAll.Range("F21").Value = 69
$wsAll.Range("F22").Value = 27
$wsAll.Range("F23").Value = 128
